$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab) to match the new date
$ws.Name = "Through 2022-07-17"

# Update the label for the July row
$ws.Range("A8").Value = "July (through 07-17)"

# Update July row values (row 8)
$ws.Range("B8").Value = 23
$ws.Range("C8").Value = 32
$ws.Range("D8").Value = 36
$ws.Range("E8").Value = 39
$ws.Range("G8").Value = 70
$ws.Range("H8").Value = 81
$ws.Range("I8").Value = 95

# Update Total row values (row 9)
$ws.Range("B9").Value = 148
$ws.Range("C9").Value = 280
$ws.Range("D9").Value = 426
$ws.Range("E9").Value = 392
$ws.Range("G9").Value = 542
$ws.Range("H9").Value = 841
$ws.Range("I9").Value = 900
